$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("H2").Value = 0.1032715332661945

$ws.Range("B3").Value = 0.123544959004684
$ws.Range("H3").Value = 0.2268164922708785

$ws.Range("B4").Value = 0.1125071070265921
$ws.Range("H4").Value = 0.2157786402927865

$ws.Range("B5").Value = 0.04584762580201045
$ws.Range("C5").Value = 0.004884665086340472
$ws.Range("D5").Value = 8.074458152572564
$ws.Range("E5").Value = 0.05013392399081189
$ws.Range("F5").Value = 0.03626765548662099
$ws.Range("G5").Value = 0.0554275961173997
$ws.Range("H5").Value = 0.1491191590682049

$ws.Range("B6").Value = 0.03483020773375735
$ws.Range("H6").Value = 0.1381017409999518

$ws.Range("B7").Value = 0.01866757594154206
$ws.Range("C7").Value = 0.002098370891433816
$ws.Range("D7").Value = 2.705064845679695
$ws.Range("E7").Value = 0.008510740238794712
$ws.Range("F7").Value = 0.01454797859067791
$ws.Range("G7").Value = 0.02278717329240586
$ws.Range("H7").Value = 0.1219391092077365

$ws.Range("B8").Value = 0.01678749391398004
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = 0.1200590271801745

$ws.Range("B9").Value = 0.01148705543886069
$ws.Range("C9").Value = 0.001817406281004231
$ws.Range("D9").Value = 2.084904913975639
$ws.Range("E9").Value = 0.008441666954877369
$ws.Range("F9").Value = 0.007919437685998357
$ws.Range("G9").Value = 0.01505467319172316
$ws.Range("H9").Value = 0.1147585887050551

$ws.Range("B10").Value = 0.01155506298580298
$ws.Range("C10").Value = 0.002082678145171691
$ws.Range("D10").Value = 2.107974717519345
$ws.Range("E10").Value = 0.009045247916819056
$ws.Range("F10").Value = 0.007467477304198136
$ws.Range("G10").Value = 0.01564264866740796
$ws.Range("H10").Value = 0.1148265962519974

$ws.Range("B11").Value = 0.02972441137369737
$ws.Range("H11").Value = 0.1329959446398918

$ws.Range("B12").Value = 0.05409485438562589
$ws.Range("H12").Value = 0.1573663876518203

$ws.Range("B13").Value = 0.06918147161944022
$ws.Range("H13").Value = 0.1724530048856347

$ws.Range("B14").Value = 0.08072206526740502
$ws.Range("H14").Value = 0.1839935985335995

$ws.Range("B15").Value = 0.08619898135964785
$ws.Range("H15").Value = 0.1894705146258423

$ws.Range("B16").Value = 0.08985738308665765
$ws.Range("H16").Value = 0.1931289163528521

$ws.Range("B17").Value = 0.09408923803573851
$ws.Range("H17").Value = 0.197360771301933

$ws.Range("B18").Value = -0.1032715332661945
$ws.Range("C18").Value = 0.008905996699974422
$ws.Range("D18").Value = -19.05966684785116
$ws.Range("E18").Value = 0.02685169165100391
$ws.Range("F18").Value = -0.120772599487219
$ws.Range("G18").Value = -0.08577046704517004
$ws.Range("H18").Value = 0

$ws.Range("B19").Value = 0.09406115305583133
$ws.Range("H19").Value = 0.1973326863220258

$ws.Range("B20").Value = 0.09819629145515409
$ws.Range("H20").Value = 0.2014678247213486

$ws.Range("B21").Value = 0.1016710335447454
$ws.Range("H21").Value = 0.2049425668109399

$ws.Range("B22").Value = 0.1058913419359381
$ws.Range("H22").Value = 0.2091628752021325

$ws.Range("B23").Value = 0.1096621026567973
$ws.Range("C23").Value = 0.007196970191805244
$ws.Range("D23").Value = 27.19819957207644
$ws.Range("E23").Value = 0.04345474295017757
$ws.Range("F23").Value = 0.09551341072141653
$ws.Range("G23").Value = 0.1238107945921782
$ws.Range("H23").Value = 0.2129336359229918

$ws.Range("B24").Value = 0.1131890454873457
$ws.Range("C24").Value = 0.007069284488701446
$ws.Range("D24").Value = 27.99589419668919
$ws.Range("E24").Value = 0.03090247037235266
$ws.Range("F24").Value = 0.09928944041638178
$ws.Range("G24").Value = 0.1270886505583096
$ws.Range("H24").Value = 0.2164605787535401

$ws.Range("B25").Value = 0.1155133960877915
$ws.Range("C25").Value = 0.007144910753670359
$ws.Range("D25").Value = 29.10440059802164
$ws.Range("E25").Value = 0.04618075995359579
$ws.Range("F25").Value = 0.1014584394361259
$ws.Range("G25").Value = 0.129568352739457
$ws.Range("H25").Value = 0.2187849293539859

$ws.Range("B26").Value = 0.1178704374144319
$ws.Range("C26").Value = 0.007252321536902761
$ws.Range("D26").Value = 29.46692268445513
$ws.Range("E26").Value = 0.0496586813913303
$ws.Range("F26").Value = 0.1036197730170177
$ws.Range("G26").Value = 0.1321211018118469
$ws.Range("H26").Value = 0.2211419706806264

$ws.Range("B27").Value = 0.1220811423598562
$ws.Range("C27").Value = 0.007146291459685908
$ws.Range("D27").Value = 29.11740688913864
$ws.Range("E27").Value = 0.05567456673404157
$ws.Range("F27").Value = 0.1080305649009582
$ws.Range("G27").Value = 0.1361317198187539
$ws.Range("H27").Value = 0.2253526756260506

$ws.Range("B28").Value = 0.1232720042135041
$ws.Range("C28").Value = 0.00716153581725948
$ws.Range("D28").Value = 26.16275451701625
$ws.Range("E28").Value = 0.08403230770418317
$ws.Range("F28").Value = 0.1091963941166053
$ws.Range("G28").Value = 0.1373476143104025
$ws.Range("H28").Value = 0.2265435374796985

$ws.Range("B29").Value = 0.01434479949641376
$ws.Range("C29").Value = 0.001418444108487877
$ws.Range("D29").Value = 2.575276824397787
$ws.Range("E29").Value = 0.00891385920699507
$ws.Range("F29").Value = 0.01156192527061761
$ws.Range("G29").Value = 0.01712767372220998
$ws.Range("H29").Value = 0.1176163327626082
